$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price snapshot refresh (GitHub Actions cron job).
# Columns D (Price) and E (Volume change) are stored as plain text in the feed:
# prices reuse "." as both thousands and decimal separator (e.g. "26.080.05"),
# so most are never valid numbers, and percentages keep their padding spaces.
# A few new price values do parse as plain numbers (e.g. "218.07"), and Excel
# would silently coerce those to the Number type on assignment, so each such
# cell is pre-formatted as Text to keep it a string like the rest of the column.

$ws.Range("D2").Value = '26.054.77'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.650.56'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.07'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5299'
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2611'
$ws.Range("E8").Value = '  -2.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06289'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.31'
$ws.Range("E10").Value = '  -3.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07745'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.472'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.651.70'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5442'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '0.0₅8101'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.14'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = '26.088.47'
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.554'
$ws.Range("E19").Value = '  -2.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.59'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.03'
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.987'
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '139.65'
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1242'
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.253'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.23'
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05914'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.276'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.501'
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.239'
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.542'
$ws.Range("E33").Value = '  -6.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.417'
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9413'
$ws.Range("E35").Value = '  -4.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.754'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5659'
$ws.Range("E37").Value = '  -4.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01604'
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.845'
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8458'
$ws.Range("E40").Value = '  -2.13%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '1.008.28'
$ws.Range("E42").Value = '  -3.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.64'
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").Value = '1.802.91'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '56.94'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9997'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4297'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.478'
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.838'
$ws.Range("E51").Value = '  -3.49%  '

Write-Host "Applied cryptos list update"
